$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fixed naive component forecaster bug - updated YoY forecast values (columns C and E, rows 2-39)
$ws.Range("C2").Value = 3.123541145015474
$ws.Range("E2").Value = 4.320516327661528
$ws.Range("C3").Value = 3.959010658874851
$ws.Range("E3").Value = 3.433494243648449
$ws.Range("C4").Value = 4.722695063536686
$ws.Range("E4").Value = 7.855477094481422
$ws.Range("C5").Value = 6.739021039846627
$ws.Range("E5").Value = -0.7259153295281151
$ws.Range("C6").Value = 2.619839412265601
$ws.Range("E6").Value = -0.17790865651377
$ws.Range("C7").Value = -0.7919564768266385
$ws.Range("E7").Value = 2.233697987812078
$ws.Range("C8").Value = 1.877689851450803
$ws.Range("E8").Value = 2.959667200710037
$ws.Range("C9").Value = 2.705004599189187
$ws.Range("E9").Value = 0.5406927319912658
$ws.Range("C10").Value = 1.110374544249249
$ws.Range("E10").Value = 3.29487077883559
$ws.Range("C11").Value = 2.267566233338814
$ws.Range("E11").Value = 3.15890982365572
$ws.Range("C12").Value = 2.688433258834588
$ws.Range("E12").Value = 3.239674285955152
$ws.Range("C13").Value = 1.014079695989589
$ws.Range("E13").Value = 2.571626871154176
$ws.Range("C14").Value = 3.013853578092252
$ws.Range("E14").Value = 2.652928973511215
$ws.Range("C15").Value = 1.331295149770684
$ws.Range("E15").Value = -0.3738725857433511
$ws.Range("C16").Value = 0.04589006555719699
$ws.Range("E16").Value = 1.021287096146906
$ws.Range("C17").Value = 0.009546395482029624
$ws.Range("E17").Value = 0.5726247744375135
$ws.Range("C18").Value = 0.8709390141433015
$ws.Range("E18").Value = 0.6379602509701376
$ws.Range("C19").Value = 0.7652063367885598
$ws.Range("E19").Value = 2.297544413125596
$ws.Range("C20").Value = 2.267579219134386
$ws.Range("E20").Value = 2.825914290412324
$ws.Range("C21").Value = 3.146753122914103
$ws.Range("E21").Value = 2.510325059131513
$ws.Range("C22").Value = 1.769033835366818
$ws.Range("E22").Value = -1.563640406432543
$ws.Range("C23").Value = -4.774715709990263
$ws.Range("E23").Value = 1.685921024959058
$ws.Range("C24").Value = 1.95493704440024
$ws.Range("E24").Value = 3.595026567604331
$ws.Range("C25").Value = 3.478075069442799
$ws.Range("E25").Value = 2.123512403013161
$ws.Range("C26").Value = 1.232342134690434
$ws.Range("E26").Value = 1.023960954496861
$ws.Range("C27").Value = 0.2542811494408159
$ws.Range("E27").Value = 1.414810393331356
$ws.Range("C28").Value = 1.519778766382096
$ws.Range("E28").Value = 0.3626330124320232
$ws.Range("C29").Value = 1.469441753880329
$ws.Range("E29").Value = 1.339091979913909
$ws.Range("C30").Value = 1.638203081492495
$ws.Range("E30").Value = 1.006355688239569
$ws.Range("C31").Value = 2.268697431234346
$ws.Range("E31").Value = 3.257368055312471
$ws.Range("C32").Value = 1.984425467899631
$ws.Range("E32").Value = -0.3299132127116078
$ws.Range("C33").Value = 0.6066448776129052
$ws.Range("E33").Value = 0.2691274977562275
$ws.Range("C34").Value = -4.243076347305386
$ws.Range("E34").Value = 23.52713615747899
$ws.Range("C35").Value = 1.438499295329754
$ws.Range("E35").Value = 6.182015844361843
$ws.Range("C36").Value = 1.906593537051537
$ws.Range("E36").Value = 1.559766133975371
$ws.Range("C37").Value = 0.08348019664223827
$ws.Range("E37").Value = -0.3513378667146627
$ws.Range("C38").Value = -0.214505326882275
$ws.Range("E38").Value = 0.256124515548195
$ws.Range("C39").Value = 0.1651547428133782
$ws.Range("E39").Value = -0.02867520550564606
